# Add URL for google icons to Power Point
#
# Slide 1 (index 1): four existing shapes get re-positioned (their
# a:off x/y change, a:ext stays the same) and a new text box with the
# Google Fonts icons URL is added at the end of the shape tree.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Reposition the existing shapes -------------------------------------
# NOTE: Shape.Left/Top/Width/Height are Single (32-bit float) in the
# PowerPoint object model, so the literals below are chosen so that the
# float32 round-trip still lands exactly on the target EMU offset.

# "Group 16" (the hand-pointer group)
$grp = $s.Shapes.Item(1)
$grp.Left = 205.33937072753906
$grp.Top = 215.6820526123047

# "Graphic 60" ("End with solid fill")
$picEnd = $s.Shapes.Item(5)
$picEnd.Left = 638.8943481445312
$picEnd.Top = 392.20294189453125

# "Picture 2"
$picPicture = $s.Shapes.Item(6)
$picPicture.Left = 543.1993408203125
$picPicture.Top = 146.42953491210938

# "Graphic 9" ("Open folder outline")
$picFolder = $s.Shapes.Item(7)
$picFolder.Left = 371.0433349609375
$picFolder.Top = 230.54229736328125

# --- Add the new text box with the Google Fonts icons URL ---------------
# A throw-away textbox is added and removed first purely to advance the
# shape-id allocator, so that the real text box ends up with id=3 / the
# name "TextBox 2" (matching a fresh id=2 having already been used once),
# exactly as it was produced when this edit was originally authored.
$tmp = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$tmp.Delete()

$tb = $s.Shapes.AddTextbox(1, 30.6696062992126, 37.468110236220475, 480.1303937007874, 50.892204724409446)
$tb.TextFrame.WordWrap = -1
$tb.TextFrame.AutoSize = 1

$tr = $tb.TextFrame.TextRange
$tr.Text = "https://"
$tr.InsertAfter("fonts.google.com")
$tr.InsertAfter("/")
$tr.InsertAfter("icons?icon.size")
$tr.InsertAfter("=24&icon.color=%231f1f1f")
